$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": add a new column BB (06-aug) next to BA (05-aug)
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell BB1 - copy the formatting (bold / centered / bordered)
# from the previous header cell BA1, then set its text.
$wsPrix.Range("BB1").Value = "06-aug"
$wsPrix.Range("BA1").Copy()
$wsPrix.Range("BB1").PasteSpecial(-4122)   # xlPasteFormats

# Hourly values for the new 06-aug column (rows 2-25)
$wsPrix.Range("BB2").Value  = 75.63
$wsPrix.Range("BB3").Value  = 40.07
$wsPrix.Range("BB4").Value  = 34.8
$wsPrix.Range("BB5").Value  = 31.1
$wsPrix.Range("BB6").Value  = 41.95
$wsPrix.Range("BB7").Value  = 47.55
$wsPrix.Range("BB8").Value  = 51.35
$wsPrix.Range("BB9").Value  = 78.89
$wsPrix.Range("BB10").Value = 83.19
$wsPrix.Range("BB11").Value = 29.71
$wsPrix.Range("BB12").Value = -0.01
$wsPrix.Range("BB13").Value = -0.76
$wsPrix.Range("BB14").Value = -0.02
$wsPrix.Range("BB15").Value = -0.02
$wsPrix.Range("BB16").Value = -1.02
$wsPrix.Range("BB17").Value = -0.03
$wsPrix.Range("BB18").Value = 3.18
$wsPrix.Range("BB19").Value = 52
$wsPrix.Range("BB20").Value = 80.31999999999999
$wsPrix.Range("BB21").Value = 106.22
$wsPrix.Range("BB22").Value = 114.64
$wsPrix.Range("BB23").Value = 117.02
$wsPrix.Range("BB24").Value = 110
$wsPrix.Range("BB25").Value = 90.54000000000001

# ---------------------------------------------------------------------
# Sheet "Gaz": append a new row (2025-08-04) after the last row (50)
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date cell to be stored as plain text (like the rest of the
# column) instead of letting Excel auto-convert it to a date serial.
$wsGaz.Range("A51").NumberFormat = "@"
$wsGaz.Range("A51").Value = "2025-08-04"
$wsGaz.Range("A51").Style = "Normal"

$wsGaz.Range("B51").Value = 33.525

# ---------------------------------------------------------------------
# Sheet "CO2": append a new row (2025-08-04) after the last row (50)
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A51").NumberFormat = "@"
$wsCo2.Range("A51").Value = "2025-08-04"
$wsCo2.Range("A51").Style = "Normal"

$wsCo2.Range("B51").Value = 70.27
